$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = 39
$row = 40

# Copy formatting from the row above for every column, then set values.
$ws.Cells.Item($srcRow, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item($row, 1).Value = 39

$ws.Cells.Item($row, 2).Value = "india"
$ws.Cells.Item($row, 3).Value = "isl"
$ws.Cells.Item($row, 4).Value = "2023-2024"

$ws.Cells.Item($srcRow, 5).Copy()
$ws.Cells.Item($row, 5).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item($row, 5).Value = 45259.64583333334

$ws.Cells.Item($row, 6).Value = "Kerala Blasters"
$ws.Cells.Item($row, 7).Value = 3
$ws.Cells.Item($row, 8).Value = "Chennaiyin"
$ws.Cells.Item($row, 9).Value = 3

$ws.Cells.Item($row, 10).Value = 2.17
$ws.Cells.Item($row, 11).Value = "25/11/2023 18:13"
$ws.Cells.Item($row, 12).Value = 1.85
$ws.Cells.Item($row, 13).Value = "29/11/2023 15:22"
$ws.Cells.Item($row, 14).Value = 3.51
$ws.Cells.Item($row, 15).Value = "25/11/2023 18:13"
$ws.Cells.Item($row, 16).Value = 3.75
$ws.Cells.Item($row, 17).Value = "29/11/2023 15:22"
$ws.Cells.Item($row, 18).Value = 3.24
$ws.Cells.Item($row, 19).Value = "25/11/2023 18:13"
$ws.Cells.Item($row, 20).Value = 4.17
$ws.Cells.Item($row, 21).Value = "29/11/2023 15:22"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/india/isl/kerala-blasters-chennaiyin-fc/x6LzikV2/"

$excel.CutCopyMode = 0
